$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PC1 (column B) and PC2 (column C) values for rows 2-17
# reflecting the generalized variable (material properties in addition to crack geometry)
$ws.Range("B2").Value = -0.09604866752742558
$ws.Range("C2").Value = -0.2051611755279131
$ws.Range("B3").Value = -0.2620242827054002
$ws.Range("C3").Value = -0.5567197273228244
$ws.Range("B4").Value = -0.2264729653632841
$ws.Range("C4").Value = -0.01064773026437103
$ws.Range("B5").Value = 0.3397273705783637
$ws.Range("C5").Value = -0.4209851171209288
$ws.Range("B6").Value = -0.5024545140173816
$ws.Range("C6").Value = -0.3191206482617719
$ws.Range("B7").Value = -0.1608068732529925
$ws.Range("C7").Value = -0.1850300299348308
$ws.Range("B8").Value = 0.1479190898505301
$ws.Range("C8").Value = -0.1652696216258207
$ws.Range("B9").Value = -0.04126112377577518
$ws.Range("C9").Value = -0.09912481737059048
$ws.Range("B10").Value = 0.2356590894368341
$ws.Range("C10").Value = 0.2131754087559259
$ws.Range("B11").Value = -0.5092295915877625
$ws.Range("C11").Value = 0.289107960977763
$ws.Range("B12").Value = -0.1685846637637685
$ws.Range("C12").Value = 0.3594981498025291
$ws.Range("B13").Value = 0.1075963397527032
$ws.Range("C13").Value = -0.1576161664434762
$ws.Range("B14").Value = -0.2738559048800113
$ws.Range("C14").Value = 0.04645526402658837
$ws.Range("B15").Value = -0.06274528781670646
$ws.Range("C15").Value = -0.002042344924774836
$ws.Range("B16").Value = 0.1319776103530128
$ws.Range("C16").Value = -0.01771674187207919
$ws.Range("B17").Value = 0.04807753920908057
$ws.Range("C17").Value = -0.1087714465602982
